$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = 45917
$ws.Range("A22").NumberFormat = "mm-dd-yy"
$ws.Range("B22").Value = 6
$ws.Range("C22").Value = 182

$ws.Range("G31").Select()
